$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 46 (Leve Item ID 4584)
$ws.Range("H46").Value = 1433266.9
$ws.Range("J46").Value = 5478
$ws.Range("L46").Value = 16434
$ws.Range("N46").Value = -16672
# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 562.6
$ws.Range("I58").Value = 562.6
$ws.Range("K58").Value = 1687.8
$ws.Range("M58").Value = -1537.8
# Row 60 (Leve Item ID 4584)
$ws.Range("H60").Value = 1433266.9
$ws.Range("J60").Value = 5478
$ws.Range("L60").Value = 16434
$ws.Range("N60").Value = -17402
# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 7617.304
$ws.Range("I116").Value = 9362.875
$ws.Range("K116").Value = 9362.875
$ws.Range("M116").Value = -5920.875
# Row 134 (Leve Item ID 41997)
$ws.Range("H134").Value = 88861.625
$ws.Range("J134").Value = 88861.625
$ws.Range("L134").Value = 88861.625
$ws.Range("N134").Value = -99001.625
# Row 136 (Leve Item ID 42164)
$ws.Range("H136").Value = 89947
$ws.Range("J136").Value = 89947
$ws.Range("L136").Value = 89947
$ws.Range("N136").Value = -100147
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 2560.158
$ws.Range("I137").Value = 2352.75
$ws.Range("J137").Value = 2915.7144
$ws.Range("K137").Value = 7058.25
$ws.Range("L137").Value = 8747.143199999999
$ws.Range("M137").Value = -4508.25
$ws.Range("N137").Value = -13847.1432

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 4586
$ws.Range("I61").Value = 4527
$ws.Range("K61").Value = 4527
$ws.Range("M61").Value = -4315
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 2387
$ws.Range("I74").Value = 2439.7727
$ws.Range("K74").Value = 2439.7727
$ws.Range("M74").Value = -1565.7727
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 2387
$ws.Range("I77").Value = 2439.7727
$ws.Range("K77").Value = 12198.8635
$ws.Range("M77").Value = -7830.863499999999
# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 4586
$ws.Range("I136").Value = 4527
$ws.Range("K136").Value = 13581
$ws.Range("M136").Value = -11031
# Row 137 (Leve Item ID 43227)
$ws.Range("H137").Value = 86566.336
$ws.Range("J137").Value = 89849.5
$ws.Range("L137").Value = 89849.5
$ws.Range("N137").Value = -100049.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 6676.7856
$ws.Range("I105").Value = 7155.5557
$ws.Range("J105").Value = 5815
$ws.Range("K105").Value = 7155.5557
$ws.Range("L105").Value = 5815
$ws.Range("M105").Value = -5408.5557
$ws.Range("N105").Value = -9309
# Row 111 (Leve Item ID 25789)
$ws.Range("H111").Value = 60000
$ws.Range("J111").Value = 60000
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -68180
# Row 132 (Leve Item ID 41855)
$ws.Range("H132").Value = 79886.25
$ws.Range("J132").Value = 79886.25
$ws.Range("L132").Value = 79886.25
$ws.Range("N132").Value = -90006.25
# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 3421.4443
$ws.Range("I134").Value = 2899.16
$ws.Range("K134").Value = 8697.48
$ws.Range("M134").Value = -6162.48

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 2579.8293
$ws.Range("I58").Value = 2625.625
$ws.Range("K58").Value = 2625.625
$ws.Range("M58").Value = -2422.625
# Row 123 (Leve Item ID 35334)
$ws.Range("H123").Value = 89999
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -99799
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2801.0667
$ws.Range("I134").Value = 1423.6923
$ws.Range("K134").Value = 4271.0769
$ws.Range("M134").Value = -1736.0769
# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 2579.8293
$ws.Range("I136").Value = 2625.625
$ws.Range("K136").Value = 7876.875
$ws.Range("M136").Value = -5326.875

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5 (Leve Item ID 43974)
$ws.Range("H5").Value = 751.75
$ws.Range("I5").Value = 644.8570999999999
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 1934.5713
$ws.Range("L5").Value = 4500
$ws.Range("M5").Value = -1822.5713
$ws.Range("N5").Value = -4724
# Row 38 (Leve Item ID 4860)
$ws.Range("H38").Value = 93.71429000000001
$ws.Range("I38").Value = 96
$ws.Range("K38").Value = 288
$ws.Range("M38").Value = 59
# Row 39 (Leve Item ID 4712)
$ws.Range("H39").Value = 6645.5864
$ws.Range("J39").Value = 6747.393
$ws.Range("L39").Value = 20242.179
$ws.Range("N39").Value = -20830.179
# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 6869
$ws.Range("J55").Value = 6869
$ws.Range("L55").Value = 20607
$ws.Range("N55").Value = -20961
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1188.2
$ws.Range("I68").Value = 647.6667
$ws.Range("K68").Value = 1943.0001
$ws.Range("M68").Value = -1132.0001
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1188.2
$ws.Range("I71").Value = 647.6667
$ws.Range("K71").Value = 5829.0003
$ws.Range("M71").Value = -1773.0003
# Row 113 (Leve Item ID 27843)
$ws.Range("H113").Value = 1206
$ws.Range("I113").Value = 450
$ws.Range("J113").Value = 1458
$ws.Range("K113").Value = 1350
$ws.Range("L113").Value = 4374
$ws.Range("M113").Value = 820
$ws.Range("N113").Value = -8714
# Row 116 (Leve Item ID 27866)
$ws.Range("H116").Value = 1600
$ws.Range("I116").Value = 1600
$ws.Range("K116").Value = 4800
$ws.Range("M116").Value = -1358
# Row 135 (Leve Item ID 43974)
$ws.Range("H135").Value = 751.75
$ws.Range("I135").Value = 644.8570999999999
$ws.Range("J135").Value = 1500
$ws.Range("K135").Value = 5803.7139
$ws.Range("L135").Value = 13500
$ws.Range("M135").Value = -3268.7139
$ws.Range("N135").Value = -18570

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 5899.3
$ws.Range("J70").Value = 6001
$ws.Range("L70").Value = 6001
$ws.Range("N70").Value = -6541
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 5899.3
$ws.Range("J73").Value = 6001
$ws.Range("L73").Value = 6001
$ws.Range("N73").Value = -7873
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 400
$ws.Range("I113").Value = 400
$ws.Range("K113").Value = 400
$ws.Range("M113").Value = 1770
# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 64993
$ws.Range("J130").Value = 64993
$ws.Range("L130").Value = 64993
$ws.Range("N130").Value = -75033
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 5542.8
$ws.Range("I132").Value = 5165.3335
$ws.Range("K132").Value = 15496.0005
$ws.Range("M132").Value = -12966.0005
# Row 135 (Leve Item ID 42006)
$ws.Range("H135").Value = 89000
$ws.Range("J135").Value = 89000
$ws.Range("L135").Value = 89000
$ws.Range("N135").Value = -99140
# Row 140 (Leve Item ID 42458)
$ws.Range("H140").Value = 79760
$ws.Range("J140").Value = 79760
$ws.Range("L140").Value = 79760
$ws.Range("N140").Value = -90120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 96 (Leve Item ID 19735)
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
# Row 100 (Leve Item ID 19995)
$ws.Range("H100").Value = 2505.1177
$ws.Range("I100").Value = 2753.6365
$ws.Range("J100").Value = 2049.5
$ws.Range("K100").Value = 2753.6365
$ws.Range("L100").Value = 2049.5
$ws.Range("M100").Value = -2212.6365
$ws.Range("N100").Value = -3131.5
# Row 111 (Leve Item ID 25820)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3937.4
$ws.Range("I122").Value = 4004.4285
$ws.Range("K122").Value = 12013.2855
$ws.Range("M122").Value = -9563.2855
# Row 123 (Leve Item ID 35408)
$ws.Range("H123").Value = 89999
$ws.Range("J123").Value = 89999
$ws.Range("L123").Value = 89999
$ws.Range("N123").Value = -99799
# Row 125 (Leve Item ID 34271)
$ws.Range("H125").Value = 84832.5
$ws.Range("J125").Value = 84832.5
$ws.Range("L125").Value = 84832.5
$ws.Range("N125").Value = -94672.5
# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3436.4
$ws.Range("I132").Value = 3450
$ws.Range("K132").Value = 10350
$ws.Range("M132").Value = -7820
# Row 134 (Leve Item ID 42024)
$ws.Range("H134").Value = 66699.2
$ws.Range("J134").Value = 66699.2
$ws.Range("L134").Value = 66699.2
$ws.Range("N134").Value = -76839.2
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 1425.7037
$ws.Range("I136").Value = 1259.56
$ws.Range("K136").Value = 3778.68
$ws.Range("M136").Value = -1228.68
# Row 138 (Leve Item ID 42334)
$ws.Range("H138").Value = 79681.60000000001
$ws.Range("J138").Value = 79681.60000000001
$ws.Range("L138").Value = 79681.60000000001
$ws.Range("N138").Value = -89961.60000000001
# Row 139 (Leve Item ID 43310)
$ws.Range("H139").Value = 89665.664
$ws.Range("J139").Value = 89665.664
$ws.Range("L139").Value = 89665.664
$ws.Range("N139").Value = -99945.664
# Row 141 (Leve Item ID 42487)
$ws.Range("H141").Value = 89999
$ws.Range("J141").Value = 89999
$ws.Range("L141").Value = 89999
$ws.Range("N141").Value = -100359

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 57 (Leve Item ID 10927)
$ws.Range("H57").Value = 101899.336
$ws.Range("J57").Value = 101899.336
$ws.Range("L57").Value = 101899.336
$ws.Range("N57").Value = -103407.336
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2099.3333
$ws.Range("I126").Value = 1899
$ws.Range("K126").Value = 5697
$ws.Range("M126").Value = -3227
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1655.098
$ws.Range("I136").Value = 929.2889
$ws.Range("J136").Value = 7098.6665
$ws.Range("K136").Value = 2787.8667
$ws.Range("L136").Value = 21295.9995
$ws.Range("M136").Value = -237.8667
$ws.Range("N136").Value = -26395.9995
# Row 138 (Leve Item ID 42347)
$ws.Range("H138").Value = 75214.5
$ws.Range("I138").Value = 60000
$ws.Range("K138").Value = 60000
$ws.Range("M138").Value = -54860
# Row 141 (Leve Item ID 42505)
$ws.Range("H141").Value = 69999
$ws.Range("I141").Value = 69999
$ws.Range("K141").Value = 69999
$ws.Range("M141").Value = -64819
